# Week 15 simulations update
# Appends newly simulated play-by-play figures to the running per-play
# lists (shared strings) and refreshes the aggregate totals on the
# YDS / OFF / DEF / ST / TURNS / PEN sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Helper: append a space separated list of numbers to the existing
# text of a cell (the cell already holds a space separated list).
# ---------------------------------------------------------------
function Append-Values {
    param(
        $Range,
        [string]$ToAppend
    )
    $cur = $Range.Text
    $Range.Value = $cur + " " + $ToAppend
}

# ---------------------------------------------------------------
# YDS sheet - rushing/passing yardage lists for OFF (B) and DEF (C)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("YDS")
Append-Values $ws.Range("B2") "7 8 4 4 2 0 7 -1 2 7 5 1 7 5 7 4 9 5 -3 13 8 5 8 6 8 7"
Append-Values $ws.Range("B3") "5 15 6 60 3 7 -2 1 17 3 9 6 18 10 1 8 24"
Append-Values $ws.Range("C2") "0 14 9 5 5 9 1 1 0 2 17 2 3 3 12 4 -1 6 4 8 2 4 2 6 1 3 2 10 9 0 3 2 0 2 4"
Append-Values $ws.Range("C3") "15 7 5 11 9 -2 15 12 10 25 7 2 7 5 17 59 7 36 5 11 8 3 1"

# ---------------------------------------------------------------
# OFF sheet - aggregate totals for RATT / PATT rows
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("OFF")
$ws.Range("C2").Value = 325
$ws.Range("D2").Value = 18
$ws.Range("F2").Value = 117
$ws.Range("G2").Value = 84
$ws.Range("I2").Value = 13
$ws.Range("J2").Value = 36
$ws.Range("N2").Value = 41
$ws.Range("O2").Value = 42

$ws.Range("C3").Value = 329
$ws.Range("E3").Value = 59
$ws.Range("F3").Value = 191
$ws.Range("G3").Value = 48
$ws.Range("H3").Value = 61
$ws.Range("I3").Value = 118
$ws.Range("J3").Value = 94
$ws.Range("L3").Value = 533
$ws.Range("M3").Value = 326
$ws.Range("Q3").Value = 925

# ---------------------------------------------------------------
# DEF sheet - aggregate totals for RATT / PATT rows
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("DEF")
$ws.Range("C2").Value = 371
$ws.Range("D2").Value = 21
$ws.Range("F2").Value = 121
$ws.Range("G2").Value = 110
$ws.Range("I2").Value = 12
$ws.Range("J2").Value = 57
$ws.Range("N2").Value = 30

$ws.Range("C3").Value = 329
$ws.Range("D3").Value = 7
$ws.Range("E3").Value = 57
$ws.Range("F3").Value = 210
$ws.Range("G3").Value = 60
$ws.Range("H3").Value = 31
$ws.Range("I3").Value = 115
$ws.Range("J3").Value = 104
$ws.Range("L3").Value = 544
$ws.Range("M3").Value = 372
$ws.Range("Q3").Value = 978

# ---------------------------------------------------------------
# ST sheet - kickoff (KO, B/C columns) and punt (PT, D column)
# counts plus the per-kick detail lists (distance / return att /
# return yards)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ST")
$ws.Range("B2").Value = 127
$ws.Range("D2").Value = 118
$ws.Range("F2").Value = 359
$ws.Range("G2").Value = 342
$ws.Range("H2").Value = 9
$ws.Range("I2").Value = 6

Append-Values $ws.Range("B4") "58 71"
Append-Values $ws.Range("B5") "0 29"
Append-Values $ws.Range("B6") "17 17"
Append-Values $ws.Range("D3") "50 46 51 63 18"
Append-Values $ws.Range("D4") "0 0 0 0 0"
Append-Values $ws.Range("D5") "8 0 0"

# ---------------------------------------------------------------
# TURNS sheet - Road turnover totals
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("TURNS")
$ws.Range("B3").Value = 13
$ws.Range("D3").Value = 15

# ---------------------------------------------------------------
# PEN sheet - False start penalty count
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("PEN")
$ws.Range("B2").Value = 30
